$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: append two new sentences (each wrapped in its own run, with
# a standalone space-run on either side) to the paragraph that ends with
# "...for the satisfaction of our client."
# ---------------------------------------------------------------------

$sentence = "Once taken and reviewed all possible risks, we will have a successful number, about how long we could take this project and how many times we would have to restart the project in case of finding something a failure, locating and focusing on the error, to maintain the quality and the final time of the product."

# Extend the existing (fully-formatted) run so the new text inherits
# its exact run properties (rFonts incl. cs, sz, szCs) via Find/Replace.
$needle = "for the satisfaction of our client."
$replacement = "for the satisfaction of our client. " + $sentence + " "
$fr = $d.Content
$found = $fr.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

# Locate the freshly-inserted sentence text so we can split it (and its
# surrounding spaces) into their own separate <w:r> runs, mirroring the
# target diff, while keeping every run's inherited formatting intact.
$locator = $d.Content
$locator.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentStart = $locator.Start
$sentEnd = $locator.End

$space1 = $d.Range($sentStart - 1, $sentStart)
$space1.Bold = 1
$space1.Bold = 0

$sentRange = $d.Range($sentStart, $sentEnd)
$sentRange.Bold = 1
$sentRange.Bold = 0

$space2 = $d.Range($sentEnd, $sentEnd + 1)
$space2.Bold = 1
$space2.Bold = 0

# ---------------------------------------------------------------------
# Change 2: insert a brand-new, empty paragraph (carrying the usual
# Times New Roman / 12pt paragraph mark formatting) right after the
# paragraph that holds the "_GoBack" bookmark, before the final
# paragraph of the document.
# ---------------------------------------------------------------------

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>')
